# Add a "Save" column (column H) to the s_vals sheet, matching the
# header style used by the other header cells (bold, bordered, centered),
# and fill the data rows with 0 values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "sum" header (G1) onto the new
# header cell so the new header picks up the same style (bold font,
# border, centered alignment) instead of Excel's default style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New "Save" data values for the existing rows.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
